$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.194.62"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "1.881.07"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.50"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5140"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3912"
$ws.Range("E8").Value = "  +2.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08352"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.56"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.236"
$ws.Range("E12").Value = "  +0.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.68"
$ws.Range("E13").Value = "  +1.13%  "

$ws.Range("D14").Value = "1.874.18"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.264"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.21"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06679"
$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.81"
$ws.Range("E20").Value = "  +0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.051"
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").Value = "28.226.41"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  +0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.259"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").Value = "2.087.83"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.491"
$ws.Range("E27").Value = "  -2.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.94"
$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.65"
$ws.Range("E29").Value = "  +1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.26"
$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1063"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.041"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.844"
$ws.Range("E33").Value = "  +4.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.634"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02453"
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06563"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2190"
$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.199"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6505"
$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.998"
$ws.Range("E41").Value = "  +2.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.225"
$ws.Range("E42").Value = "  -1.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.34"
$ws.Range("E43").Value = "  +0.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6147"
$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.02"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.286"
$ws.Range("E46").Value = "  +0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.677"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.019"
$ws.Range("E48").Value = "  +1.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.232"
$ws.Range("E49").Value = "  +2.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.77"
$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06915"
$ws.Range("E51").Value = "  +0.88%  "
